$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Batch upload" terminology cleanup: the two reviewer-initial columns move
# after NPM in the shared-string table and are renamed Penguji -> Reviewer.
# The actual header cells (A1=NPM, B1=Tanggal Seminar, C1=Link Konferensi)
# keep their text; only D1/E1 text changes.
$ws.Range("D1").Value = "Inisial Dosen Reviewer 1"
$ws.Range("E1").Value = "Inisial Dosen Reviewer 2"

# Column widths were resized (A and C widened considerably, D/E nudged).
$ws.Columns.Item(1).ColumnWidth = 18.333333333333336
$ws.Columns.Item(3).ColumnWidth = 29.0
$ws.Columns.Item(4).ColumnWidth = 20.0
$ws.Columns.Item(5).ColumnWidth = 22.333333333333336

# Cursor/selection left on H7 when the sheet was last saved.
$ws.Range("H7").Select() | Out-Null
